$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 508.5
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 508.5
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H34").Value = 10889.6
$ws.Range("J34").Value = 20049
$ws.Range("L34").Value = 20049
$ws.Range("N34").Value = -20455
$ws.Range("H36").Value = 10889.6
$ws.Range("J36").Value = 20049
$ws.Range("L36").Value = 20049
$ws.Range("N36").Value = -21479
$ws.Range("H113").Value = 17775
$ws.Range("I113").Value = 15746.5
$ws.Range("K113").Value = 15746.5
$ws.Range("M113").Value = -12492.5
$ws.Range("H135").Value = 4236.7144
$ws.Range("J135").Value = 16220.75
$ws.Range("L135").Value = 145986.75
$ws.Range("N135").Value = -151056.75
$ws.Range("H137").Value = 10245.854
$ws.Range("I137").Value = 2589.3809
$ws.Range("J137").Value = 18285.15
$ws.Range("K137").Value = 7768.1427
$ws.Range("L137").Value = 54855.45
$ws.Range("M137").Value = -5218.1427
$ws.Range("N137").Value = -59955.45
$ws.Range("H138").Value = 2541.92
$ws.Range("I138").Value = 2291
$ws.Range("J138").Value = 2918.3
$ws.Range("K138").Value = 6873
$ws.Range("L138").Value = 8754.900000000001
$ws.Range("M138").Value = -1733
$ws.Range("N138").Value = -19034.9
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14217.308
$ws.Range("I32").Value = 4904.1924
$ws.Range("K32").Value = 4904.1924
$ws.Range("M32").Value = -4617.1924
$ws.Range("H61").Value = 28967
$ws.Range("I61").Value = 37870.668
$ws.Range("J61").Value = 26538.727
$ws.Range("K61").Value = 37870.668
$ws.Range("L61").Value = 26538.727
$ws.Range("M61").Value = -37658.668
$ws.Range("N61").Value = -26962.727
$ws.Range("H102").Value = 12646.823
$ws.Range("I102").Value = 937.3125
$ws.Range("K102").Value = 937.3125
$ws.Range("M102").Value = 684.6875
$ws.Range("H132").Value = 2641679.8
$ws.Range("I132").Value = 4942.4736
$ws.Range("K132").Value = 14827.4208
$ws.Range("M132").Value = -12297.4208
$ws.Range("H136").Value = 28967
$ws.Range("I136").Value = 37870.668
$ws.Range("J136").Value = 26538.727
$ws.Range("K136").Value = 113612.004
$ws.Range("L136").Value = 79616.181
$ws.Range("M136").Value = -111062.004
$ws.Range("N136").Value = -84716.181
$ws.Range("H139").Value = 65483.25
$ws.Range("I139").Value = 49000
$ws.Range("J139").Value = 66981.73
$ws.Range("K139").Value = 49000
$ws.Range("L139").Value = 66981.73
$ws.Range("M139").Value = -43860
$ws.Range("N139").Value = -77261.73
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 954.5
$ws.Range("I99").Value = 818.13336
$ws.Range("K99").Value = 818.13336
$ws.Range("M99").Value = 679.86664
$ws.Range("H132").Value = 96246.875
$ws.Range("J132").Value = 96246.875
$ws.Range("L132").Value = 96246.875
$ws.Range("N132").Value = -106366.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 23477.25
$ws.Range("J18").Value = 23477.25
$ws.Range("L18").Value = 23477.25
$ws.Range("N18").Value = -23937.25
$ws.Range("H31").Value = 14606.333
$ws.Range("I31").Value = 5469.125
$ws.Range("K31").Value = 5469.125
$ws.Range("M31").Value = -5174.125
$ws.Range("H34").Value = 14606.333
$ws.Range("I34").Value = 5469.125
$ws.Range("K34").Value = 5469.125
$ws.Range("M34").Value = -5267.125
$ws.Range("H82").Value = 65111
$ws.Range("J82").Value = 65111
$ws.Range("L82").Value = 65111
$ws.Range("N82").Value = -65833
$ws.Range("H85").Value = 65111
$ws.Range("J85").Value = 65111
$ws.Range("L85").Value = 65111
$ws.Range("N85").Value = -67607
$ws.Range("H94").Value = 1701.6
$ws.Range("I94").Value = 2254.5
$ws.Range("K94").Value = 2254.5
$ws.Range("M94").Value = -1803.5
$ws.Range("H100").Value = 66237.25
$ws.Range("J100").Value = 66237.25
$ws.Range("L100").Value = 66237.25
$ws.Range("N100").Value = -68401.25
$ws.Range("H122").Value = 5241.794
$ws.Range("I122").Value = 3530.1667
$ws.Range("K122").Value = 10590.5001
$ws.Range("M122").Value = -8140.500100000001
$ws.Range("H132").Value = 6526.433
$ws.Range("I132").Value = 2275.7368
$ws.Range("J132").Value = 13868.546
$ws.Range("K132").Value = 6827.2104
$ws.Range("L132").Value = 41605.638
$ws.Range("M132").Value = -4297.2104
$ws.Range("N132").Value = -46665.638
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1627464.6
$ws.Range("J5").Value = 3050356
$ws.Range("L5").Value = 9151068
$ws.Range("N5").Value = -9151292
$ws.Range("H135").Value = 1627464.6
$ws.Range("J135").Value = 3050356
$ws.Range("L135").Value = 27453204
$ws.Range("N135").Value = -27458274
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9132.857
$ws.Range("I113").Value = 9132.857
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 9132.857
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6962.857
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2887.6897
$ws.Range("I122").Value = 2369.158
$ws.Range("J122").Value = 3872.9
$ws.Range("K122").Value = 7107.474
$ws.Range("L122").Value = 11618.7
$ws.Range("M122").Value = -4657.474
$ws.Range("N122").Value = -16518.7
$ws.Range("H126").Value = 4706.393
$ws.Range("J126").Value = 5638
$ws.Range("L126").Value = 16914
$ws.Range("N126").Value = -21854
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 5576.303
$ws.Range("I132").Value = 2566.8572
$ws.Range("K132").Value = 7700.571599999999
$ws.Range("M132").Value = -5170.571599999999
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 75958
$ws.Range("J140").Value = 75958
$ws.Range("L140").Value = 75958
$ws.Range("N140").Value = -86318
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6475.222
$ws.Range("J40").Value = 8017.273
$ws.Range("L40").Value = 8017.273
$ws.Range("N40").Value = -8289.273000000001
$ws.Range("H61").Value = 3167.2
$ws.Range("I61").Value = 1778.5883
$ws.Range("J61").Value = 6118
$ws.Range("K61").Value = 1778.5883
$ws.Range("L61").Value = 6118
$ws.Range("M61").Value = -1576.5883
$ws.Range("N61").Value = -6522
$ws.Range("H68").Value = 10632.143
$ws.Range("I68").Value = 7481.25
$ws.Range("J68").Value = 14833.333
$ws.Range("K68").Value = 7481.25
$ws.Range("L68").Value = 14833.333
$ws.Range("M68").Value = -6732.25
$ws.Range("N68").Value = -16331.333
$ws.Range("H71").Value = 10632.143
$ws.Range("I71").Value = 7481.25
$ws.Range("J71").Value = 14833.333
$ws.Range("K71").Value = 37406.25
$ws.Range("L71").Value = 74166.66500000001
$ws.Range("M71").Value = -33662.25
$ws.Range("N71").Value = -81654.66500000001
$ws.Range("H113").Value = 3167.2
$ws.Range("I113").Value = 1778.5883
$ws.Range("J113").Value = 6118
$ws.Range("K113").Value = 1778.5883
$ws.Range("L113").Value = 6118
$ws.Range("M113").Value = 391.4117000000001
$ws.Range("N113").Value = -10458
$ws.Range("H132").Value = 982568.1
$ws.Range("I132").Value = 3625.6086
$ws.Range("J132").Value = 2233439
$ws.Range("K132").Value = 10876.8258
$ws.Range("L132").Value = 6700317
$ws.Range("M132").Value = -8346.825800000001
$ws.Range("N132").Value = -6705377
$ws.Range("H136").Value = 12673.907
$ws.Range("I136").Value = 12697.3
$ws.Range("J136").Value = 12653.565
$ws.Range("K136").Value = 38091.89999999999
$ws.Range("L136").Value = 37960.695
$ws.Range("M136").Value = -35541.89999999999
$ws.Range("N136").Value = -43060.695
$ws.Range("H140").Value = 149441.9
$ws.Range("J140").Value = 149441.9
$ws.Range("L140").Value = 149441.9
$ws.Range("N140").Value = -159801.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 19555.62
$ws.Range("I54").Value = 18603.5
$ws.Range("K54").Value = 18603.5
$ws.Range("M54").Value = -18083.5
$ws.Range("H96").Value = 2804.9
$ws.Range("I96").Value = 2166.6667
$ws.Range("K96").Value = 2166.6667
$ws.Range("M96").Value = -793.6667000000002
$ws.Range("H100").Value = 865.3333
$ws.Range("I100").Value = 800.3333
$ws.Range("J100").Value = 897.8333
$ws.Range("K100").Value = 1600.6666
$ws.Range("L100").Value = 1795.6666
$ws.Range("M100").Value = -1059.6666
$ws.Range("N100").Value = -2877.6666
$ws.Range("H132").Value = 6907.6665
$ws.Range("I132").Value = 2230
$ws.Range("K132").Value = 6690
$ws.Range("M132").Value = -4160
$ws.Range("H136").Value = 6546.868
$ws.Range("I136").Value = 1734.0571
$ws.Range("J136").Value = 15905.111
$ws.Range("K136").Value = 5202.1713
$ws.Range("L136").Value = 47715.333
$ws.Range("M136").Value = -2652.1713
$ws.Range("N136").Value = -52815.333
